# Revert the two "CU - 23" / "CU - 24" rows (Consultar profesores / Consultar
# clientes) that had been merged into the "Casos de Uso" sheet, restoring the
# sheet to its pre-merge state: rows 27 and 28 go back to being blank (the
# formatting/styles stay, only the values are removed), and the selected
# cell reverts to C23 (its pre-merge location) instead of C28.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Clear out the data that the merge had added in B27:I28, leaving the
# underlying cell styles untouched.
$ws.Range("B27:I28").ClearContents()

# Restore the previously active cell/selection on this sheet.
$ws.Range("C23").Select()
